$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 188, shifting existing rows 188:271 down to 189:272
$ws.Rows.Item(188).Insert()

# Populate the newly inserted row 188 with the new record's data
$ws.Cells.Item(188, 1).Value = 5
$ws.Cells.Item(188, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(188, 3).Value = "Maule"
$ws.Cells.Item(188, 4).Value = 44523
$ws.Cells.Item(188, 5).Value = 7
$ws.Cells.Item(188, 6).Value = 100112043
$ws.Cells.Item(188, 7).Value = "Pepino ensalada"
$ws.Cells.Item(188, 8).Value = "Sin especificar"
$ws.Cells.Item(188, 9).Value = "Primera"
$ws.Cells.Item(188, 10).Value = 400
$ws.Cells.Item(188, 11).Value = 8000
$ws.Cells.Item(188, 12).Value = 8000
$ws.Cells.Item(188, 13).Value = 8000
$ws.Cells.Item(188, 14).Value = "`$/caja 80 unidades"
$ws.Cells.Item(188, 15).Value = "Región del Maule"
$ws.Cells.Item(188, 16).Value = 100
$ws.Cells.Item(188, 17).Value = 80
$ws.Cells.Item(188, 18).Value = "Hortaliza"
